# Add a new "Cancer_Related_Project" column (N) to the cost-recovery sheet,
# fix the mis-placed "tissue type" value (was sitting alone in row 9's M
# cell; it belongs to row 8, Dr. Sanford Barsky's record), and backfill a
# couple of missing numeric values that had landed in the wrong row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix data that had drifted into the wrong cells -------------------------
# G4 (276) actually belongs to row 5 (Dr. Dash Chandravanu), not row 4.
$ws.Range("G4").ClearContents()
$ws.Range("G5").Value = 276
$ws.Range("I5").Value = 276

# Row 8 (Dr. Sanford Barsky) should show "Breast cancer cell line" as its
# tissue type instead of "mouse"; the stray value previously lived by
# itself down in M9.
$ws.Range("M8").Value = "Breast cancer cell line"
$ws.Range("M9").ClearContents()

# --- Add the new "Cancer_Related_Project" column (N) ------------------------
$ws.Range("N2").Value = "no"
$ws.Range("N3").Value = "no"
$ws.Range("N4").Value = "yes"
$ws.Range("N5").Value = "no"
$ws.Range("N6").Value = "no"
$ws.Range("N7").Value = "no"
$ws.Range("N8").Value = "yes"

$ws.Range("N1").Value = "Cancer_Related_Project"
$ws.Range("N1").Font.Bold = $true

# --- Update the selection to match where the editor left off ---------------
[void]$ws.Range("N1").Select()
